# Apply the changes described by the diff:
# 1. Update the "Date" metadata value on the "Metadata" sheet.
# 2. Update the canonical terminology URLs on the "Elements" sheet.
# 3. Widen column Z on the "Elements" sheet to fit the new (longer) URLs.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# 1. Update the Date value (B8) on the Metadata sheet
$wsMetadata.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# 2. Update the canonical terminology URLs on the Elements sheet
$wsElements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R53-DiplomePaysEEE/FHIR/TRE-R53-DiplomePaysEEE?vs"

# 3. Widen column Z (26th column) to fit the new (longer) URL text
$wsElements.Columns.Item(26).ColumnWidth = 73.8
